$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 21: updated dct:modified timestamp ---
$ws.Range("B21").Value = '2023-09-22T14:37:15+00:00'

# --- Fix row 128: typo correction + vocab id correction ---
$ws.Range("B128").Value = 'whole air matrix'
$ws.Range("G128").Value = 'vocab:1098'

# --- Fix rows 148-154: vocab.1101 -> vocab:1101 ---
foreach ($r in 148..154) {
    $ws.Range("G$r").Value = 'vocab:1101'
}

# --- Append new rows 175-196 (new matrix vocabulary terms) ---
$newRows = @(
    ,@('vocab:1151', 'blood-whole blood matrix', 'BWB', 'Blood-whole blood', 'vocab:1172')
    ,@('vocab:1152', 'blood -plasma matrix', 'BP', 'Blood -plasma', 'vocab:1172')
    ,@('vocab:1153', 'blood -serum matrix', 'BS', 'Blood -serum', 'vocab:1172')
    ,@('vocab:1154', 'cord blood-whole blood matrix', 'CBWB', 'Cord blood-whole blood', 'vocab:1172')
    ,@('vocab:1155', 'cord blood-plasma matrix', 'CBP', 'Cord blood-plasma', 'vocab:1172')
    ,@('vocab:1156', 'cord blood-serum matrix', 'CBS', 'Cord blood-serum', 'vocab:1172')
    ,@('vocab:1157', 'urine-spot matrix', 'US', 'Urine-spot', 'vocab:1172')
    ,@('vocab:1158', 'urine-24h matrix', 'UD', 'Urine-24h', 'vocab:1172')
    ,@('vocab:1159', 'urine-morning urine matrix', 'UM', 'Urine-morning urine', 'vocab:1172')
    ,@('vocab:1160', 'saliva and/or sputum matrix', 'SA', 'Saliva and or sputum', 'vocab:1172')
    ,@('vocab:1161', 'semen matrix', 'SEM', 'Semen', 'vocab:1172')
    ,@('vocab:1162', 'hair matrix', 'H', 'Hair', 'vocab:1172')
    ,@('vocab:1163', 'exhaled breath condensate matrix', 'EBC', 'Exhaled breath condensate', 'vocab:1172')
    ,@('vocab:1164', 'red blood cells matrix', 'RBC', 'Red blood cells', 'vocab:1172')
    ,@('vocab:1165', 'breast milk matrix', 'BM', 'Breast milk', 'vocab:1172')
    ,@('vocab:1166', 'adipose tissue/fat matrix', 'ADI', 'Adipose tissue/fat', 'vocab:1172')
    ,@('vocab:1167', 'all toe nails matrix', 'ATN', 'All toe nails', 'vocab:1172')
    ,@('vocab:1168', 'big toe nails matrix', 'BTN', 'Big toe nails', 'vocab:1172')
    ,@('vocab:1169', 'dermal wipes matrix', 'DW', 'Dermal wipes', 'vocab:1172')
    ,@('vocab:1170', 'amniotic fluid matrix', 'AF', 'Amniotic fluid', 'vocab:1172')
    ,@('vocab:1171', 'placenta tissue matrix', 'PLT', 'Placenta tissue', 'vocab:1172')
    ,@('vocab:1172', 'human matrix', '', '', '')
)

$startRow = 175
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = $row[4]
    for ($c = 8; $c -le 42; $c++) {
        $ws.Cells.Item($r, $c).Value = ""
    }
}
